$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = -5.874999999999999
$ws.Range("B9").Value = 8.562400000000006
$ws.Range("D11").Value = -8.383999999999997
$ws.Range("B18").Value = 4.966900000000003
$ws.Range("B20").Value = 5.473899999999998
$ws.Range("E21").Value = 13.2965
